# "adding DSI to Experience"
#
# The "General Assembly: 13 week immersive course ..." sentence in the
# Experience section gets "13 week" turned into "13-week": the space
# between "13" and "week" becomes a hyphen, the now-unnecessary
# gramStart/gramEnd proofing-error bracket around "13 week" is removed,
# and the floating "_GoBack" bookmark (left over from the author's last
# edit position) is relocated to sit right after the new hyphen.

$d = $word.ActiveDocument

# --- locate the paragraph's text span --------------------------------
$rFirst = $d.Content
$rFirst.Find.Execute("General Assembly: ", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$paraStart = $rFirst.Start

$rGap = $d.Content
$rGap.Find.Execute("13 week", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$gapPos = $rGap.Start + 2          # the space between "13" and "week"

$rLast = $d.Content
$rLast.Find.Execute("Eagle Ford Shale Oil and Gas Play. ", $false, $false, `
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraEnd = $rLast.End

# --- rebuild the whole sentence as explicit runs ----------------------
# (replacing the full span lets us drop the gramStart/gramEnd proofErr
# pair that bracketed "13 week", since they sit just inside this range)
$newRuns = (
    '<w:r><w:t>General Assembly: 13</w:t></w:r>' +
    '<w:r><w:t>-</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">week immersive course applying data ' +
    'collection, cleaning, exploratory analysis, modelling, data ' +
    'visualization, and various machine learning </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>tools</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> to solve real-world data science ' +
    'problems. Created a Na' + [char]0x00EF + 've Bayes and Random Forest ' +
    'Classifier to discern between 10,000+ Reddit posts using Natural ' +
    'Language Processing techniques. Trained a Convolutional Neural ' +
    'Network to assist in diagnosing COVID-19 pneumonia presentation in ' +
    'chest x-rays. Capstone project utilized forecasting methods to ' +
    'better understand mechanism driving induced earthquakes in the ' +
    'Eagle Ford Shale Oil and Gas Play. </w:t></w:r>'
)

$xmlFrag = (
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" ' +
    'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $newRuns + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
)

$rWhole = $d.Range($paraStart, $paraEnd)
$rWhole.InsertXML($xmlFrag)

# --- relocate the hidden "_GoBack" bookmark to sit right after the
#     freshly-typed hyphen, matching where Word leaves it after an edit.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
$rBookmark = $d.Range($gapPos + 1, $gapPos + 1)
$d.Bookmarks.Add("_GoBack", $rBookmark) | Out-Null

Write-Output "done"
